# Updates Leve profit/price figures across sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 4163.636
$ws.Range("I64").Value = 4400
$ws.Range("J64").Value = 4111.1113
$ws.Range("K64").Value = 4400
$ws.Range("L64").Value = 4111.1113
$ws.Range("M64").Value = -4152
$ws.Range("N64").Value = -4607.1113
# Row 67
$ws.Range("H67").Value = 4163.636
$ws.Range("I67").Value = 4400
$ws.Range("J67").Value = 4111.1113
$ws.Range("K67").Value = 4400
$ws.Range("L67").Value = 4111.1113
$ws.Range("M67").Value = -3542
$ws.Range("N67").Value = -5827.1113
# Row 74
$ws.Range("H74").Value = 5249
$ws.Range("I74").Value = 2245
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 2245
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -1309
$ws.Range("N74").Value = -7872
# Row 77
$ws.Range("H77").Value = 5249
$ws.Range("I77").Value = 2245
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 11225
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -6545
$ws.Range("N77").Value = -39360
# Row 101
$ws.Range("M101").Value = 872
$ws.Range("H101").Value = 1478.3334
$ws.Range("I101").Value = 250
$ws.Range("J101").Value = 2092.5
$ws.Range("K101").Value = 750
$ws.Range("L101").Value = 6277.5
$ws.Range("N101").Value = -9521.5
# Row 113
$ws.Range("H113").Value = 33336742
$ws.Range("I113").Value = 71431576
$ws.Range("J113").Value = 3759.1875
$ws.Range("K113").Value = 71431576
$ws.Range("L113").Value = 3759.1875
$ws.Range("M113").Value = -71428322
$ws.Range("N113").Value = -10267.1875
# Row 132
$ws.Range("H132").Value = 2990.6333
$ws.Range("I132").Value = 3154.25
$ws.Range("K132").Value = 9462.75
$ws.Range("M132").Value = -6932.75
# Row 137
$ws.Range("H137").Value = 1431.6765
$ws.Range("I137").Value = 1288.8572
$ws.Range("K137").Value = 3866.5716
$ws.Range("M137").Value = -1316.5716

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1619.0588
$ws.Range("I2").Value = 1508.4286
$ws.Range("K2").Value = 1508.4286
$ws.Range("M2").Value = -1395.4286
# Row 32
$ws.Range("H32").Value = 5776.769
$ws.Range("I32").Value = 4117.84
$ws.Range("K32").Value = 4117.84
$ws.Range("M32").Value = -3830.84
# Row 88
$ws.Range("M88").ClearContents()
$ws.Range("H88").Value = 336684.66
$ws.Range("I88").Value = 0
$ws.Range("K88").Value = 0
# Row 91
$ws.Range("M91").ClearContents()
$ws.Range("H91").Value = 336684.66
$ws.Range("I91").Value = 0
$ws.Range("K91").Value = 0
# Row 109
$ws.Range("H109").Value = 19999.5
$ws.Range("J109").Value = 19999.5
$ws.Range("L109").Value = 19999.5
$ws.Range("N109").Value = -22773.5
# Row 112
$ws.Range("H112").Value = 37999
$ws.Range("J112").Value = 37999
$ws.Range("L112").Value = 37999
$ws.Range("N112").Value = -40953
# Row 116
$ws.Range("H116").Value = 1619.0588
$ws.Range("I116").Value = 1508.4286
$ws.Range("K116").Value = 1508.4286
$ws.Range("M116").Value = 785.5714

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1619.0588
$ws.Range("I3").Value = 1508.4286
$ws.Range("K3").Value = 1508.4286
$ws.Range("M3").Value = -1394.4286
# Row 86
$ws.Range("H86").Value = 1786.1
$ws.Range("I86").Value = 1452.5238
$ws.Range("J86").Value = 2564.4443
$ws.Range("K86").Value = 1452.5238
$ws.Range("L86").Value = 2564.4443
$ws.Range("M86").Value = -329.5237999999999
$ws.Range("N86").Value = -4810.4443
# Row 89
$ws.Range("H89").Value = 1786.1
$ws.Range("I89").Value = 1452.5238
$ws.Range("J89").Value = 2564.4443
$ws.Range("K89").Value = 7262.619
$ws.Range("L89").Value = 12822.2215
$ws.Range("M89").Value = -1646.619
$ws.Range("N89").Value = -24054.2215
# Row 99
$ws.Range("H99").Value = 1501.5714
$ws.Range("I99").Value = 1188.1333
$ws.Range("K99").Value = 1188.1333
$ws.Range("M99").Value = 309.8667
# Row 107
$ws.Range("H107").Value = 1754.625
$ws.Range("I107").Value = 1002.2
$ws.Range("J107").Value = 3008.6667
$ws.Range("K107").Value = 1002.2
$ws.Range("L107").Value = 3008.6667
$ws.Range("M107").Value = 917.8
$ws.Range("N107").Value = -6848.6667
# Row 110
$ws.Range("H110").Value = 45699.5
$ws.Range("J110").Value = 45699.5
$ws.Range("L110").Value = 45699.5
$ws.Range("N110").Value = -53879.5

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1388.5714
$ws.Range("I16").Value = 1388.5714
$ws.Range("K16").Value = 1388.5714
$ws.Range("M16").Value = -1101.5714
# Row 106
$ws.Range("N106").ClearContents()
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
# Row 113
$ws.Range("H113").Value = 1388.5714
$ws.Range("I113").Value = 1388.5714
$ws.Range("K113").Value = 1388.5714
$ws.Range("M113").Value = 781.4286

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 824.8570999999999
$ws.Range("I5").Value = 845.6667
$ws.Range("K5").Value = 2537.0001
$ws.Range("M5").Value = -2425.0001
# Row 86
$ws.Range("H86").Value = 62500336
$ws.Range("I86").Value = 483.33334
$ws.Range("J86").Value = 100000250
$ws.Range("K86").Value = 1450.00002
$ws.Range("L86").Value = 300000750
$ws.Range("M86").Value = -264.0000199999999
$ws.Range("N86").Value = -300003122
# Row 89
$ws.Range("H89").Value = 62500336
$ws.Range("I89").Value = 483.33334
$ws.Range("J89").Value = 100000250
$ws.Range("K89").Value = 4350.00006
$ws.Range("L89").Value = 900002250
$ws.Range("M89").Value = 1577.99994
$ws.Range("N89").Value = -900014106
# Row 131
$ws.Range("H131").Value = 770.5
$ws.Range("J131").Value = 770.5
$ws.Range("L131").Value = 2311.5
$ws.Range("N131").Value = -12391.5
# Row 135
$ws.Range("H135").Value = 824.8570999999999
$ws.Range("I135").Value = 845.6667
$ws.Range("K135").Value = 7611.0003
$ws.Range("M135").Value = -5076.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 148149170
$ws.Range("I122").Value = 47619930
$ws.Range("K122").Value = 142859790
$ws.Range("M122").Value = -142857340
# Row 132
$ws.Range("H132").Value = 26493.088
$ws.Range("I132").Value = 5756.2
$ws.Range("K132").Value = 17268.6
$ws.Range("M132").Value = -14738.6

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 5545.364
$ws.Range("I40").Value = 3554.182
$ws.Range("J40").Value = 7536.5454
$ws.Range("K40").Value = 3554.182
$ws.Range("L40").Value = 7536.5454
$ws.Range("M40").Value = -3418.182
$ws.Range("N40").Value = -7808.5454
# Row 122
$ws.Range("H122").Value = 983162.6
$ws.Range("I122").Value = 1510734.5
$ws.Range("K122").Value = 4532203.5
$ws.Range("M122").Value = -4529753.5
# Row 136
$ws.Range("H136").Value = 28738.555
$ws.Range("I136").Value = 42541.168
$ws.Range("K136").Value = 127623.504
$ws.Range("M136").Value = -125073.504

$ws = $wb.Worksheets.Item("WVR")
# Row 104
$ws.Range("N104").ClearContents()
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
# Row 126
$ws.Range("N126").ClearContents()
$ws.Range("H126").Value = 698.4
$ws.Range("I126").Value = 698.4
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2095.2
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 374.8000000000002

Write-Output "Applied Typhon_Profits updates"
